$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 1")
Write-Host $ws.Name
